# AutoCommit_26 июня 2024 г. 12:15:10_SibNout2023
#
# Adds the "key" row (R32:U32 = 2,3,4,5 — same pattern as R3:U3) and a
# totals row (R34:U34 = SUM of each column over rows 4:30) below the
# existing student table, then updates the view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 32: repeat of the key values found in row 3 (R3:U3) ---
$ws.Range("R32").Value = 2
$ws.Range("S32").Value = 3
$ws.Range("T32").Value = 4
$ws.Range("U32").Value = 5

# --- New row 34: column totals over the student rows (4:30) ---
$ws.Range("R34").Formula = "=SUM(R4:R30)"
$ws.Range("S34:U34").FormulaR1C1 = "=SUM(R[-30]C:R[-4]C)"

# --- View state: keep the existing freeze split (2 cols / 3 rows),
#     scroll the frozen pane so H13 is the first visible cell, and move
#     the active selection to V27 ---
$win = $excel.ActiveWindow
$excel.Goto($ws.Range("H13"), $true) | Out-Null
$ws.Range("V27").Select() | Out-Null
